$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 28739
$ws.Range("E2").Value = 1475
$ws.Range("F2").Value = 1475
$ws.Range("G2").Value = 1128
$ws.Range("H2").Value = 943
$ws.Range("I2").Value = 659
$ws.Range("J2").Value = 284
$ws.Range("K2").Value = 27732
$ws.Range("L2").Value = 19461
$ws.Range("M2").Value = 8270
$ws.Range("N2").Value = 5593
$ws.Range("O2").Value = 2677
$ws.Range("P2").Value = 662
$ws.Range("Q2").Value = 2197
$ws.Range("R2").Value = -1936
$ws.Range("S2").Value = 28
$ws.Range("T2").Value = 1735
$ws.Range("U2").Value = 463
$ws.Range("V2").Value = 10503
$ws.Range("W2").Value = 5.13
$ws.Range("X2").Value = 3.28
$ws.Range("Y2").Value = 13.36
$ws.Range("Z2").Value = 3.66
$ws.Range("AA2").Value = 235.32
$ws.Range("AB2").Value = 750.88
$ws.Range("AC2").Value = 5254
$ws.Range("AD2").Value = 15.27
$ws.Range("AE2").Value = 43374
$ws.Range("AF2").Value = 1.85
$ws.Range("AG2").Value = 350
$ws.Range("AH2").Value = 0.44
$ws.Range("AI2").Value = 6.85
$ws.Range("AJ2").Value = 13247561

# Row 3
$ws.Range("D3").Value = 28071
$ws.Range("E3").Value = 1063
$ws.Range("F3").Value = 1063
$ws.Range("G3").Value = -199
$ws.Range("H3").Value = -374
$ws.Range("I3").Value = -430
$ws.Range("J3").Value = 56
$ws.Range("K3").Value = 24955
$ws.Range("L3").Value = 15939
$ws.Range("M3").Value = 9016
$ws.Range("N3").Value = 5609
$ws.Range("O3").Value = 3407
$ws.Range("P3").Value = 662
$ws.Range("Q3").Value = 1412
$ws.Range("R3").Value = 1434
$ws.Range("S3").Value = -1464
$ws.Range("T3").Value = 1125
$ws.Range("U3").Value = 288
$ws.Range("V3").Value = 7729
$ws.Range("W3").Value = 3.79
$ws.Range("X3").Value = -1.33
$ws.Range("Y3").Value = -7.68
$ws.Range("Z3").Value = -1.42
$ws.Range("AA3").Value = 176.79
$ws.Range("AB3").Value = 751.97
$ws.Range("AC3").Value = -3246
$ws.Range("AD3").Value = -19.47
$ws.Range("AE3").Value = 43498
$ws.Range("AF3").Value = 1.45
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 0.79
$ws.Range("AI3").Value = -15
$ws.Range("AJ3").Value = 13247561

# Row 4
$ws.Range("D4").Value = 29218
$ws.Range("E4").Value = 2137
$ws.Range("F4").Value = 2137
$ws.Range("G4").Value = 1941
$ws.Range("H4").Value = 1333
$ws.Range("I4").Value = 720
$ws.Range("J4").Value = 613
$ws.Range("K4").Value = 27419
$ws.Range("L4").Value = 17305
$ws.Range("M4").Value = 10114
$ws.Range("N4").Value = 6204
$ws.Range("O4").Value = 3910
$ws.Range("P4").Value = 662
$ws.Range("Q4").Value = 3183
$ws.Range("R4").Value = -1286
$ws.Range("S4").Value = -743
$ws.Range("T4").Value = 731
$ws.Range("U4").Value = 2452
$ws.Range("V4").Value = 7540
$ws.Range("W4").Value = 7.31
$ws.Range("X4").Value = 4.56
$ws.Range("Y4").Value = 12.19
$ws.Range("Z4").Value = 5.09
$ws.Range("AA4").Value = 171.1
$ws.Range("AB4").Value = 841.77
$ws.Range("AC4").Value = 5436
$ws.Range("AD4").Value = 10.3
$ws.Range("AE4").Value = 48112
$ws.Range("AF4").Value = 1.16
$ws.Range("AG4").Value = 550
$ws.Range("AH4").Value = 0.98
$ws.Range("AI4").Value = 9.85
$ws.Range("AJ4").Value = 13247561

# Row 5
$ws.Range("D5").Value = 33925
$ws.Range("E5").Value = 2656
$ws.Range("F5").Value = 2656
$ws.Range("G5").Value = 2407
$ws.Range("H5").Value = 1854
$ws.Range("I5").Value = 996
$ws.Range("J5").Value = 858
$ws.Range("K5").Value = 29448
$ws.Range("L5").Value = 17927
$ws.Range("M5").Value = 11521
$ws.Range("N5").Value = 7015
$ws.Range("O5").Value = 4505
$ws.Range("P5").Value = 662
$ws.Range("Q5").Value = 2809
$ws.Range("R5").Value = -2925
$ws.Range("S5").Value = -393
$ws.Range("T5").Value = 2347
$ws.Range("U5").Value = 462
$ws.Range("V5").Value = 7601
$ws.Range("W5").Value = 7.83
$ws.Range("X5").Value = 5.46
$ws.Range("Y5").Value = 15.07
$ws.Range("Z5").Value = 6.52
$ws.Range("AA5").Value = 155.61
$ws.Range("AB5").Value = 970.14
$ws.Range("AC5").Value = 7522
$ws.Range("AD5").Value = 9.11
$ws.Range("AE5").Value = 54400
$ws.Range("AF5").Value = 1.26
$ws.Range("AG5").Value = 650
$ws.Range("AH5").Value = 0.95
$ws.Range("AI5").Value = 8.41
$ws.Range("AJ5").Value = 13247561

# Row 6
$ws.Range("D6").Value = 37112
$ws.Range("E6").Value = 2697
$ws.Range("F6").Value = 2697
$ws.Range("G6").Value = 2698
$ws.Range("H6").Value = 2066
$ws.Range("I6").Value = 1069
$ws.Range("K6").Value = 33979
$ws.Range("L6").Value = 19240
$ws.Range("M6").Value = 14739
$ws.Range("N6").Value = 8452
$ws.Range("P6").Value = 662
$ws.Range("Q6").Value = 3173
$ws.Range("R6").Value = -2255
$ws.Range("S6").Value = 900
$ws.Range("T6").Value = 3467
$ws.Range("U6").Value = -294
$ws.Range("V6").Value = 7806
$ws.Range("W6").Value = 7.27
$ws.Range("X6").Value = 5.57
$ws.Range("Y6").Value = 13.82
$ws.Range("Z6").Value = 6.51
$ws.Range("AA6").Value = 130.54
$ws.Range("AB6").Value = 1179.2
$ws.Range("AC6").Value = 8070
$ws.Range("AD6").Value = 6.69
$ws.Range("AE6").Value = 64785
$ws.Range("AF6").Value = 0.83
$ws.Range("AG6").Value = 750
$ws.Range("AH6").Value = 1.39
$ws.Range("AI6").Value = 9.15
$ws.Range("AJ6").Value = 13247561

# Clear rows 7-9 (columns D through AI) - erroneous estimate data removed
$ws.Range("D7:AI9").ClearContents()

Write-Output "done"